# Build site at 2022-09-26 16:07:08 UTC
# Rework of the LOM3088 course sheet: several label/value cells were
# reshuffled and the two trailing "Requisito fraco" rows were removed
# (the sheet now ends at row 24 instead of row 26).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 10: Objetivos value becomes the "Clodoaldo Saron" string ---
$ws.Range("B10").Value = "5840897 - Clodoaldo Saron"
$ws.Range("C10").Value = "5840897 - Clodoaldo Saron"

# --- Row 13: now holds "Programa resumido:" label + the 15/07/2016 value ---
$ws.Range("A13").Value = "Programa resumido:"
$ws.Range("B13").Value = "15/07/2016"
$ws.Range("C13").Value = "15/07/2016"
$ws.Rows.Item(13).RowHeight = 60

# --- Row 14: now holds "Short syllabus:" label only ---
$ws.Range("A14").Value = "Short syllabus:"
$ws.Range("B14").ClearContents()
$ws.Range("C14").ClearContents()
$ws.Rows.Item(14).RowHeight = 60

# --- Row 15: label becomes "Programa:" ; value becomes Clodoaldo Saron ---
$ws.Range("A15").Value = "Programa:"
$ws.Range("B15").Value = "5840897 - Clodoaldo Saron"
$ws.Range("C15").Value = "5840897 - Clodoaldo Saron"
$ws.Rows.Item(15).RowHeight = 120

# --- Row 16: label becomes "Syllabus:" ---
$ws.Range("A16").Value = "Syllabus:"
$ws.Rows.Item(16).RowHeight = 120

# --- Row 17: label becomes "Avaliação:", no B/C, default row height ---
$ws.Range("A17").Value = "Avaliação:"
$ws.Rows.Item(17).AutoFit()

# --- Row 18: label becomes "Método:" ; value becomes Fábio Herbst Florenzano ---
$ws.Range("A18").Value = "Método:"
$ws.Range("B18").Value = "1033242 - Fábio Herbst Florenzano"
$ws.Range("C18").Value = "1033242 - Fábio Herbst Florenzano"
$ws.Rows.Item(18).RowHeight = 60

# --- Row 19: label becomes "Critério:" ; gains the "Experimentos..." text ---
$ws.Range("A19").Value = "Critério:"
$ws.Range("B19").Value = "Experimentos desenvolvidos em laboratório didático; realização de relatórios para cada experimento."
$ws.Range("C19").Value = "Experimentos desenvolvidos em laboratório didático; realização de relatórios para cada experimento."
$ws.Rows.Item(19).RowHeight = 60

# --- Row 20: label becomes "Norma de recuperação:" ; value becomes the "Média aritmética..." text ---
$ws.Range("A20").Value = "Norma de recuperação:"
$ws.Range("B20").Value = "Média aritmética das notas obtidas nos relatórios. Será aprovado o aluno que obtiver nota final maior ou igual a 5,0."
$ws.Range("C20").Value = "Média aritmética das notas obtidas nos relatórios. Será aprovado o aluno que obtiver nota final maior ou igual a 5,0."

# --- Row 21: label becomes "Bibliografia:" ; value becomes the "Devido às características..." text ---
$ws.Range("A21").Value = "Bibliografia:"
$ws.Range("B21").Value = "Devido às características práticas da disciplina, não será oferecida recuperação."
$ws.Range("C21").Value = "Devido às características práticas da disciplina, não será oferecida recuperação."
$ws.Rows.Item(21).RowHeight = 120

# --- Row 22: label becomes "Requisitos:", no B/C, default row height ---
$ws.Range("A22").Value = "Requisitos:"
$ws.Range("B22").ClearContents()
$ws.Range("C22").ClearContents()
$ws.Rows.Item(22).AutoFit()

# --- Row 23: drop the "Bibliografia:" label, value becomes the first requisito fraco ---
$ws.Range("A23").ClearContents()
$ws.Range("B23").Value = "LOM3057 -  Introdução aos Materiais Poliméricos  (Requisito fraco)" + [char]10
$ws.Range("C23").Value = "LOM3057 -  Introdução aos Materiais Poliméricos  (Requisito fraco)" + [char]10
$ws.Rows.Item(23).RowHeight = 30

# --- Row 24: drop the "Requisitos:" label, value becomes the second requisito fraco ---
$ws.Range("A24").ClearContents()
$ws.Range("B24").Value = "LOM3058 -  Química de Polímeros  (Requisito fraco)" + [char]10
$ws.Range("C24").Value = "LOM3058 -  Química de Polímeros  (Requisito fraco)" + [char]10
$ws.Rows.Item(24).RowHeight = 30

# --- Rows 25-26 no longer exist; remove them entirely (shrinks the used range) ---
$ws.Rows("25:26").Delete()
